$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 348.6
$ws.Range("I33").Value = 238.375
$ws.Range("J33").Value = 474.57144
$ws.Range("K33").Value = 238.375
$ws.Range("L33").Value = 474.57144
$ws.Range("M33").Value = -9.375
$ws.Range("N33").Value = -932.5714399999999
$ws.Range("H107").Value = 471.4
$ws.Range("I107").Value = 229
$ws.Range("K107").Value = 229
$ws.Range("M107").Value = 1691
$ws.Range("H118").Value = 2733
$ws.Range("I118").Value = 1200
$ws.Range("K118").Value = 3600
$ws.Range("M118").Value = -1943
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 13373
$ws.Range("J44").Value = 13373
$ws.Range("L44").Value = 13373
$ws.Range("N44").Value = -14349
$ws.Range("H61").Value = 6498.5
$ws.Range("I61").Value = 5164.1665
$ws.Range("J61").Value = 8500
$ws.Range("K61").Value = 5164.1665
$ws.Range("L61").Value = 8500
$ws.Range("M61").Value = -4952.1665
$ws.Range("N61").Value = -8924
$ws.Range("H92").Value = 54999.668
$ws.Range("J92").Value = 54999.668
$ws.Range("L92").Value = 54999.668
$ws.Range("N92").Value = -59991.668
$ws.Range("H122").Value = 2032.619
$ws.Range("I122").Value = 1566.0625
$ws.Range("K122").Value = 4698.1875
$ws.Range("M122").Value = -2248.1875
$ws.Range("H136").Value = 6498.5
$ws.Range("I136").Value = 5164.1665
$ws.Range("J136").Value = 8500
$ws.Range("K136").Value = 15492.4995
$ws.Range("L136").Value = 25500
$ws.Range("M136").Value = -12942.4995
$ws.Range("N136").Value = -30600
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 940.125
$ws.Range("I80").Value = 574.4286
$ws.Range("K80").Value = 574.4286
$ws.Range("M80").Value = 423.5714
$ws.Range("H83").Value = 940.125
$ws.Range("I83").Value = 574.4286
$ws.Range("K83").Value = 2872.143
$ws.Range("M83").Value = 2119.857
$ws.Range("H125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -59840
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 97932
$ws.Range("J51").Value = 97932
$ws.Range("L51").Value = 97932
$ws.Range("N51").Value = -99404
$ws.Range("H61").Value = 97932
$ws.Range("J61").Value = 97932
$ws.Range("L61").Value = 97932
$ws.Range("N61").Value = -98628
$ws.Range("H88").Value = 17499.5
$ws.Range("J88").Value = 17499.5
$ws.Range("L88").Value = 17499.5
$ws.Range("N88").Value = -18311.5
$ws.Range("H91").Value = 17499.5
$ws.Range("J91").Value = 17499.5
$ws.Range("L91").Value = 17499.5
$ws.Range("N91").Value = -20307.5
$ws.Range("H99").Value = 2296.2856
$ws.Range("I99").Value = 2272.5
$ws.Range("K99").Value = 2272.5
$ws.Range("M99").Value = -774.5
$ws.Range("H108").Value = 34684
$ws.Range("J108").Value = 34684
$ws.Range("L108").Value = 34684
$ws.Range("N108").Value = -42364
$ws.Range("H110").Value = 50000
$ws.Range("J110").Value = 50000
$ws.Range("L110").Value = 50000
$ws.Range("N110").Value = -58180
$ws.Range("H126").Value = 2296.2856
$ws.Range("I126").Value = 2272.5
$ws.Range("K126").Value = 6817.5
$ws.Range("M126").Value = -4347.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 161.66667
$ws.Range("I12").Value = 121.4
$ws.Range("J12").Value = 181.8
$ws.Range("K12").Value = 364.2
$ws.Range("L12").Value = 545.4000000000001
$ws.Range("M12").Value = -191.2
$ws.Range("N12").Value = -891.4000000000001
$ws.Range("H14").Value = 462.58334
$ws.Range("I14").Value = 462.58334
$ws.Range("K14").Value = 1387.75002
$ws.Range("M14").Value = -1214.75002
$ws.Range("H21").Value = 375
$ws.Range("I21").Value = 375
$ws.Range("K21").Value = 1125
$ws.Range("M21").Value = -952
$ws.Range("H131").Value = 2543.4375
$ws.Range("J131").Value = 2543.4375
$ws.Range("L131").Value = 7630.3125
$ws.Range("N131").Value = -17710.3125
$ws.Range("H140").Value = 2177.7856
$ws.Range("I140").Value = 1874.0834
$ws.Range("K140").Value = 5622.2502
$ws.Range("M140").Value = -442.2502000000004
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1330.7693
$ws.Range("I102").Value = 1330.7693
$ws.Range("K102").Value = 1330.7693
$ws.Range("M102").Value = 291.2307000000001
$ws.Range("H113").Value = 4997.5
$ws.Range("I113").Value = 4997.5
$ws.Range("K113").Value = 4997.5
$ws.Range("M113").Value = -2827.5
$ws.Range("H126").Value = 2749.75
$ws.Range("I126").Value = 2749.75
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 8249.25
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5779.25
$ws.Range("N126").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 4000.8
$ws.Range("I3").Value = 4001.3333
$ws.Range("J3").Value = 4000
$ws.Range("K3").Value = 4001.3333
$ws.Range("L3").Value = 4000
$ws.Range("M3").Value = -3889.3333
$ws.Range("N3").Value = -4224
$ws.Range("H7").Value = 6158.6665
$ws.Range("I7").Value = 3998
$ws.Range("J7").Value = 6590.8
$ws.Range("K7").Value = 3998
$ws.Range("L7").Value = 6590.8
$ws.Range("M7").Value = -3886
$ws.Range("N7").Value = -6814.8
$ws.Range("H15").Value = 4000.8
$ws.Range("I15").Value = 4001.3333
$ws.Range("J15").Value = 4000
$ws.Range("K15").Value = 4001.3333
$ws.Range("L15").Value = 4000
$ws.Range("M15").Value = -3831.3333
$ws.Range("N15").Value = -4340
$ws.Range("H43").Value = 10000
$ws.Range("J43").Value = 10000
$ws.Range("L43").Value = 10000
$ws.Range("N43").Value = -10386
$ws.Range("H126").Value = 6158.6665
$ws.Range("I126").Value = 3998
$ws.Range("J126").Value = 6590.8
$ws.Range("K126").Value = 11994
$ws.Range("L126").Value = 19772.4
$ws.Range("M126").Value = -9524
$ws.Range("N126").Value = -24712.4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 90000
$ws.Range("J46").Value = 90000
$ws.Range("L46").Value = 90000
$ws.Range("N46").Value = -90462
$ws.Range("H126").Value = 4683.1904
$ws.Range("I126").Value = 2726.7693
$ws.Range("K126").Value = 8180.3079
$ws.Range("M126").Value = -5710.3079
$ws.Range("H134").Value = 90000
$ws.Range("J134").Value = 90000
$ws.Range("L134").Value = 270000
$ws.Range("N134").Value = -275070
$ws.Range("H136").Value = 3072.7878
$ws.Range("I136").Value = 2083.6365
$ws.Range("K136").Value = 6250.9095
$ws.Range("M136").Value = -3700.9095
